# Update gh-pages output: refresh "想去人数" (want-to-go headcount) figures
# in the F column of the "展览" and "全部类型" sheets, per the latest scrape.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 236
$wsExhibit.Range("F3").Value = 1081
$wsExhibit.Range("F4").Value = 72
$wsExhibit.Range("F5").Value = 405
$wsExhibit.Range("F6").Value = 74
$wsExhibit.Range("F7").Value = 540
$wsExhibit.Range("F8").Value = 58
$wsExhibit.Range("F9").Value = 6716
$wsExhibit.Range("F10").Value = 145
$wsExhibit.Range("F12").Value = 145
$wsExhibit.Range("F13").Value = 174
$wsExhibit.Range("F15").Value = 1075
$wsExhibit.Range("F16").Value = 16058
$wsExhibit.Range("F17").Value = 1577
$wsExhibit.Range("F19").Value = 325
$wsExhibit.Range("F20").Value = 173
$wsExhibit.Range("F21").Value = 115
$wsExhibit.Range("F22").Value = 11285
$wsExhibit.Range("F23").Value = 3
$wsExhibit.Range("F24").Value = 870
$wsExhibit.Range("F25").Value = 4435
$wsExhibit.Range("F26").Value = 295
$wsExhibit.Range("F27").Value = 387
$wsExhibit.Range("F28").Value = 39
$wsExhibit.Range("F29").Value = 26
$wsExhibit.Range("F32").Value = 5215

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 236
$wsAllTypes.Range("F3").Value = 1081
$wsAllTypes.Range("F4").Value = 72
$wsAllTypes.Range("F5").Value = 405
$wsAllTypes.Range("F6").Value = 74
$wsAllTypes.Range("F7").Value = 540
$wsAllTypes.Range("F9").Value = 58
$wsAllTypes.Range("F10").Value = 6716
$wsAllTypes.Range("F11").Value = 145
$wsAllTypes.Range("F13").Value = 145
$wsAllTypes.Range("F14").Value = 174
$wsAllTypes.Range("F17").Value = 1075
$wsAllTypes.Range("F18").Value = 16058
$wsAllTypes.Range("F19").Value = 1577
$wsAllTypes.Range("F21").Value = 325
$wsAllTypes.Range("F22").Value = 173
$wsAllTypes.Range("F23").Value = 115
$wsAllTypes.Range("F26").Value = 11285
$wsAllTypes.Range("F27").Value = 3
$wsAllTypes.Range("F28").Value = 870
$wsAllTypes.Range("F29").Value = 4435
$wsAllTypes.Range("F30").Value = 295
$wsAllTypes.Range("F31").Value = 387
$wsAllTypes.Range("F32").Value = 39
$wsAllTypes.Range("F33").Value = 26
$wsAllTypes.Range("F36").Value = 5215
